$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.601.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -2.16%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'1.791.28"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -2.08%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'231.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.58%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'0.5875"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -2.46%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("E7").Value = "'  +0.13%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'0.2764"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -1.14%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'23.11"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -1.87%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06723"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -4.61%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").Value = "'0.07527"
$ws.Range("D11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'1.794.99"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -2.02%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").Value = "'4.791"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.04%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'0.6116"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -2.65%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("E15").Value = "'  -2.00%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'75.21"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -4.85%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'0.000008869"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -10.26%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'28.583.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -2.29%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("D19").Value = "'5.418"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -7.15%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("E20").Value = "'  +0.07%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'209.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -6.46%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").Value = "'11.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -2.07%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("E23").Value = "'  -2.79%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +0.23%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'152.79"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -2.27%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").Value = "'8.072"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +1.26%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").Value = "'0.1257"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -3.53%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").Value = "'16.36"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -1.57%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").Value = "'1.412"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -4.63%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").Value = "'0.06146"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -4.91%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("E31").Value = "'  -1.81%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").Value = "'3.802"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.26%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("D33").Value = "'3.778"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -1.53%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").Value = "'1.732"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +0.39%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("D35").Value = "'1.045"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -5.62%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36").Value = "'0.6393"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -1.02%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").Value = "'2.502"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -1.64%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").Value = "'2.710"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -1.09%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'6.402"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -2.29%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").Value = "'0.01689"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -3.22%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").Value = "'1.140.14"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -6.06%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'0.8775"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -2.04%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("E43").Value = "'  +0.37%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = "'99.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -0.35%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'1.945.21"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.63%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").Value = "'59.78"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -4.51%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("E47").Value = "'  -4.52%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").Value = "'1.581"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +0.17%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.362"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -2.61%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05479"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -0.30%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("D51").Value = "'0.4478"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.63%  "
$ws.Range("E51").ClearFormats()
